# Update Name of Algo
# Apply updated RandomForest-imputed values to the result data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.1526
$ws.Range("A3").Value = -21.88429999999999
$ws.Range("D3").Value = -7.324999999999992
$ws.Range("D12").Value = -7.2879
$ws.Range("A14").Value = -21.77670000000001
$ws.Range("A16").Value = -21.80119999999999
$ws.Range("C18").Value = -12.6302
$ws.Range("A21").Value = -19.98859999999998
$ws.Range("A23").Value = -20.18759999999998
$ws.Range("C24").Value = -12.64239999999999
$ws.Range("D24").Value = -7.950499999999995
$ws.Range("A25").Value = -22.05629999999998
$ws.Range("C25").Value = -12.9495
$ws.Range("D25").Value = -8.652099999999995
$ws.Range("A26").Value = -21.25539999999997
$ws.Range("C27").Value = -13.2749
$ws.Range("A29").Value = -20.81679999999998
$ws.Range("C30").Value = -12.50469999999999
$ws.Range("C31").Value = -13.1527
$ws.Range("C39").Value = -12.5743
$ws.Range("A40").Value = -19.95399999999999
$ws.Range("D41").Value = -8.057999999999998
$ws.Range("C42").Value = -12.5927
$ws.Range("C48").Value = -11.5922
$ws.Range("D50").Value = -8.314400000000004
$ws.Range("C51").Value = -11.5709
$ws.Range("C52").Value = -11.2076
$ws.Range("A53").Value = -22.65810000000001
$ws.Range("D53").Value = -6.101400000000001
$ws.Range("C55").Value = -13.49539999999999
$ws.Range("C56").Value = -11.9847
$ws.Range("D56").Value = -8.019100000000003
$ws.Range("A57").Value = -21.93059999999999
$ws.Range("C57").Value = -13.03339999999999
$ws.Range("D57").Value = -8.779999999999999
$ws.Range("D58").Value = -8.219200000000003
$ws.Range("A59").Value = -22.39639999999999
$ws.Range("C60").Value = -13.09179999999999
$ws.Range("D61").Value = -7.906799999999996
$ws.Range("D63").Value = -7.934400000000005
$ws.Range("D64").Value = -7.9284
$ws.Range("A65").Value = -21.83639999999998
$ws.Range("A69").Value = -21.5707
$ws.Range("D70").Value = -7.537599999999995
$ws.Range("D72").Value = -7.415899999999999
$ws.Range("C73").Value = -12.87520000000001
$ws.Range("C74").Value = -12.00680000000001
$ws.Range("A79").Value = -20.44900000000002
$ws.Range("A83").Value = -21.9151
$ws.Range("D86").Value = -7.594399999999998
$ws.Range("C89").Value = -10.3101
$ws.Range("D89").Value = -5.627300000000003
$ws.Range("C90").Value = -12.414
$ws.Range("A91").Value = -21.38470000000002
$ws.Range("C92").Value = -10.6412
$ws.Range("A93").Value = -20.73839999999999
$ws.Range("D98").Value = -8.694599999999998
$ws.Range("A100").Value = -21.83149999999999
$ws.Range("D100").Value = -8.224600000000001
$ws.Range("D102").Value = -7.845699999999998
